$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E2" = 117
    "F2" = 86
    "H2" = 92
    "F4" = 32
    "H4" = 44
    "E5" = 174
    "F5" = 124
    "H5" = 135
    "F6" = 38
    "H6" = 48
    "E7" = 47
    "F7" = 31
    "H7" = 35
    "F8" = 7
    "H8" = 10
    "F9" = 9
    "H9" = 10
    "E10" = 753
    "F10" = 447
    "H10" = 542
    "E11" = 505
    "F11" = 311
    "H11" = 376
    "E12" = 777
    "F12" = 498
    "H12" = 584
    "E13" = 177
    "F13" = 105
    "H13" = 139
    "F14" = 91
    "H14" = 125
    "E15" = 217
    "F15" = 104
    "H15" = 155
    "E16" = 243
    "F16" = 142
    "H16" = 190
    "E17" = 133
    "F17" = 74
    "H17" = 98
    "F18" = 38
    "H18" = 55
    "E20" = 105
    "F20" = 51
    "H20" = 88
    "F21" = 94
    "H21" = 125
    "E22" = 206
    "F22" = 119
    "H22" = 161
    "F23" = 126
    "H23" = 178
    "E24" = 296
    "F24" = 175
    "H24" = 205
    "E25" = 355
    "F25" = 205
    "H25" = 265
    "E26" = 228
    "F26" = 145
    "H26" = 170
    "E27" = 411
    "F27" = 241
    "H27" = 323
    "E28" = 238
    "F28" = 117
    "H28" = 169
    "E29" = 203
    "F29" = 126
    "H29" = 167
    "E30" = 273
    "F30" = 175
    "H30" = 228
    "E31" = 86
    "F32" = 151
    "H32" = 189
    "E33" = 354
    "F33" = 191
    "H33" = 282
    "E34" = 271
    "F34" = 194
    "H34" = 232
    "F35" = 138
    "H35" = 165
    "F36" = 61
    "H36" = 71
    "E37" = 205
    "F37" = 117
    "H37" = 153
    "E38" = 112
    "F38" = 73
    "H38" = 90
    "F39" = 111
    "H39" = 162
    "E40" = 321
    "F40" = 176
    "H40" = 256
    "E41" = 461
    "F41" = 245
    "H41" = 337
    "E42" = 492
    "F42" = 305
    "H42" = 366
    "E43" = 154
    "F43" = 93
    "H43" = 120
    "E44" = 397
    "F44" = 226
    "H44" = 294
    "E45" = 192
    "F45" = 118
    "H45" = 157
    "F46" = 248
    "H46" = 312
    "E47" = 575
    "F47" = 335
    "H47" = 427
    "E48" = 291
    "F48" = 149
    "H48" = 193
    "E49" = 351
    "F49" = 182
    "H49" = 269
    "E50" = 294
    "F50" = 170
    "H50" = 243
    "F51" = 142
    "H51" = 216
    "F52" = 16
    "H52" = 24
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
